# REVUB parameters_simulation.xlsx - add HPP_name_data_lateral_flow parameter row
# (manual update: ramping envelopes capability + extended cascade calculation details)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above the current row 9 (HPP_name_data_outflow_prescribed),
# pushing all following rows down by one.
$ws.Rows(9).Insert()

# New parameter name in column A
$ws.Range("A9").Value = "HPP_name_data_lateral_flow"

# New explanatory (rich-text) description in column B
$full = ' [leave empty if unsure - only in case of cascade calculations] name of corresponding worksheet in the "data" Excel sheets (use the "inflow" workbook) from which to pull data on lateral inflow to a downstream cascade plant (hourly flow in m^3/s; hours in rows, years in columns)'
$ws.Range("B9").Value = $full

# Run 1: bold + italic lead-in note
$run1 = $ws.Range("B9").Characters(1, 63)
$run1.Font.Bold = $true
$run1.Font.Italic = $true

# Run 3: underline the word "lateral"
$run3 = $ws.Range("B9").Characters(179, 7)
$run3.Font.Underline = $true

# Row height grows to fit the longer wrapped text (4 lines instead of 3)
$ws.Rows(9).RowHeight = 58

# Reflect the cursor position left behind by the edit
$ws.Range("D9").Select()
